$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "No pic" marker in column D for the rows that got a picture
# placeholder note (matches the rows previously highlighted with style 5
# plus the final stretch of the advent calendar, per the commit).
$rows = @(4, 8, 11, 16, 19, 21, 22, 23, 24, 25)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "No pic"
}

# Move the active selection to D4 (matches the saved sheet view state).
[void]$ws.Range("D4").Select()
